$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1: add version 1.2 row (converting the old numeric 1.1 cell into the
# text value "1.2"), append a new row for the 1.2 feature bullet, and grow
# the table to match.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(3,1).NumberFormat = "@"
$ws1.Cells.Item(3,1).Value = "1.2"

$ws1.Cells.Item(4,1).Value = 1.3
$ws1.Cells.Item(4,1).NumberFormat = "@"
$ws1.Cells.Item(4,2).Value = "~ add the option to print the envelope and/or the letter"

$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:C4"))

$ws1.Range("B28").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet2: new worksheet holding the 1.3 release notes, inserted right after
# Sheet1 (so it becomes the active / second tab).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

$ws2.Columns.Item(1).ColumnWidth = 9.15
$ws2.Columns.Item(2).ColumnWidth = 83
$ws2.Columns.Item(3).ColumnWidth = 11.5

# Cell values are written in the same order the original authoring session
# created them in, so shared-string indices line up with the source file.
$ws2.Cells.Item(3,2).Value = "add a print button for the tool to go over all of the generated documents and print them"

$ws2.Cells.Item(1,1).Value = "Version"
$ws2.Cells.Item(1,2).Value = "Details"
$ws2.Cells.Item(1,3).Value = "Date"

$ws2.Cells.Item(2,1).Value = "[1.3]"
$ws2.Rows.Item(2).RowHeight = 18.75
$ws2.Cells.Item(2,2).WrapText = $true
$ws2.Cells.Item(2,2).Value = "~ print in the middle of the envelop Address1 in its first line and Address2 in the second line" + [char]10 + "~ give the user the option to generate for each contact: the letter only or envelop only or both"

$rng2 = $ws2.Range("A1:C3")
$tbl2 = $ws2.ListObjects.Add(1, $rng2, $null, 1)
$tbl2.Name = "Table2"
$tbl2.TableStyle = "TableStyleMedium4"

$ws2.Range("A2:C3").Select() | Out-Null
